# Apply the mock.xlsx data corrections (crossover detection fix in advCtrl)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("mock")

# Row 2 (Homo sapiens) - smpl3 (H)
$ws.Range("H2").Value = 1800

# Row 11 (Lactobacillus) - ctrl3 (E) and smpl3 (H)
$ws.Range("E11").Value = 6
$ws.Range("H11").Value = 601

# Row 13 (Salmonella enterica subsp. enterica) - ctrl3 (E)
$ws.Range("E13").Value = 4

# Row 14 (Root) - smpl3 (H)
$ws.Range("H14").Value = 470

# Move the active cell / selection to H15 as in the committed workbook
$ws.Range("H15").Select() | Out-Null
